$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.641.10'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '3.072.60'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '575.28'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = '169.28'
$ws.Range('E6').Value = '  -1.60%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '3.069.78'
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('E11').Value = '  -1.56%  '
$ws.Range('E12').Value = '  -3.14%  '
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('D14').Value = '35.61'
$ws.Range('E14').Value = '  -3.86%  '
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('D16').Value = '3.583.78'
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('D17').Value = '66.625.97'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Value = '16.99'
$ws.Range('E18').Value = '  +4.32%  '
$ws.Range('D19').Value = '6.96'
$ws.Range('E19').Value = '  -2.98%  '
$ws.Range('D20').Value = '3.072.72'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('D21').Value = '487.58'
$ws.Range('E21').Value = '  +2.50%  '
$ws.Range('E22').Value = '  -2.31%  '
$ws.Range('E23').Value = '  -3.88%  '
$ws.Range('D24').Value = '82.62'
$ws.Range('E24').Value = '  -1.76%  '
$ws.Range('D25').Value = '12.66'
$ws.Range('E26').Value = '  -3.15%  '
$ws.Range('D27').Value = '10.14'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('E30').Value = '  -4.93%  '
$ws.Range('E31').Value = '  -2.47%  '
$ws.Range('D32').Value = '27.55'
$ws.Range('E32').Value = '  -3.52%  '
$ws.Range('D33').Value = '0.111'
$ws.Range('E33').Value = '  -4.02%  '
$ws.Range('D34').Value = '0.0₃0908'
$ws.Range('E34').Value = '  -3.43%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '5.59'
$ws.Range('E36').Value = '  -4.59%  '
$ws.Range('D37').Value = '0.945'
$ws.Range('E37').Value = '  -3.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '47.10'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('E40').Value = '  -4.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.300'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.13%  '
$ws.Range('D42').Value = '8.28'
$ws.Range('E42').Value = '  -4.89%  '
$ws.Range('D43').Value = '2.759.38'
$ws.Range('E43').Value = '  -2.48%  '
$ws.Range('D44').Value = '2.52'
$ws.Range('E44').Value = '  -2.80%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  -3.17%  '
$ws.Range('D47').Value = '367.16'
$ws.Range('E47').Value = '  -5.01%  '
$ws.Range('D49').Value = '24.62'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('E50').Value = '  -2.01%  '
$ws.Range('E51').Value = '  -2.05%  '
